# "Generate Report for Handoff"
# The localization job for 826a9881-356d-4b87-a714-03132be8cea0.md has moved
# from "In Translation" to "Ready for handoff" (priority mt, new handoff
# timestamps). Update the Overview rollup sheet and the per-locale (zh-cn,
# de-de) detail sheets to match.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-17 14:11:42"

# ---- zh-cn detail sheet ---------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-17 14:11:37"

# ---- de-de detail sheet ----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-17 14:11:42"

# ---- Cosmetic: the Status columns got wider since "Ready for handoff" is
# longer than "In Translation" -- re-fit the affected columns so the sheets
# look right, same as Excel/the report generator would on save.
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$dede.Columns.Item(3).AutoFit()
